$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.457.84"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").Value = "1.921.64"
$ws.Range("E3").Value = "  +1.64%  "
$ws.Range("D4").Value = "'0.9991"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'243.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").Value = "'0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "'0.4696"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.65%  "
$ws.Range("D8").Value = "'0.2884"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.27%  "
$ws.Range("D9").Value = "'0.06810"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.23%  "
$ws.Range("D10").Value = "'110.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.70%  "
$ws.Range("D11").Value = "'18.43"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.83%  "
$ws.Range("D12").Value = "'0.07737"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.07%  "
$ws.Range("D13").Value = "1.892.79"
$ws.Range("D14").Value = "'5.302"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.65%  "
$ws.Range("D15").Value = "'0.6600"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("D16").Value = "'294.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.65%  "
$ws.Range("D17").Value = "30.453.77"
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("D18").Value = "'0.000007630"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "'12.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.52%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "2.140.54"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").Value = "'0.9984"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").Value = "'5.249"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.62%  "
$ws.Range("D24").Value = "'6.203"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("D25").Value = "'9.404"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("D26").Value = "'21.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.68%  "
$ws.Range("D27").Value = "'168.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("D28").Value = "'2.127"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.52%  "
$ws.Range("D29").Value = "'0.1070"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.82%  "
$ws.Range("D30").Value = "'1.366"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.06%  "
$ws.Range("D31").Value = "'4.193"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("D32").Value = "'4.003"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").Value = "'0.05044"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.13%  "
$ws.Range("D34").Value = "'0.7402"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.60%  "
$ws.Range("D35").Value = "'1.156"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.04%  "
$ws.Range("D36").Value = "'0.02096"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.02%  "
$ws.Range("D37").Value = "'2.740"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("D38").Value = "'2.676"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.22%  "
$ws.Range("D39").Value = "'2.065"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("D40").Value = "'110.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.56%  "
$ws.Range("D41").Value = "'0.8718"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.86%  "
$ws.Range("D42").Value = "'5.873"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.76%  "
$ws.Range("D43").Value = "'0.4256"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.99%  "
$ws.Range("D44").Value = "'1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D47").Value = "'7.214"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.47%  "
$ws.Range("D48").Value = "'9.306"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.18%  "
$ws.Range("D49").Value = "'0.1220"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("D50").Value = "'35.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("D51").Value = "'0.2503"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +11.99%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'67.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D46").Value = "'51.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +17.95%  "
